$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.1001513154359418
$ws.Range("J2").Value = 0.117124939437743
$ws.Range("M2").Value = 0.1727555
$ws.Range("N2").Value = 0.345511
$ws.Range("O2").Value = 0.1055017969849737
$ws.Range("P2").Value = 0.08301621587853829
$ws.Range("Q2").Value = 0.0321434065965
$ws.Range("R2").Value = 0.192860439579
$ws.Range("S2").Value = 0.01056614374890079
$ws.Range("T2").Value = 0.0097232692571244
$ws.Range("I3").Value = 0.1001513154359418
$ws.Range("J3").Value = 0.117124939437743
$ws.Range("O3").Value = 0.2846509696390457
$ws.Range("P3").Value = 0.3359750310550052
$ws.Range("S3").Value = 0.02850816904946675
$ws.Range("T3").Value = 0.03935105516491132
$ws.Range("I4").Value = 0.1001513154359418
$ws.Range("J4").Value = 0.117124939437743
$ws.Range("M4").Value = 0.268891
$ws.Range("N4").Value = 0.806673
$ws.Range("O4").Value = 0.1642117541443634
$ws.Range("P4").Value = 0.1938199938971208
$ws.Range("Q4").Value = 0.050030666133
$ws.Range("R4").Value = 0.450275995197
$ws.Range("S4").Value = 0.01644602318760146
$ws.Range("T4").Value = 0.022701155047024
$ws.Range("I5").Value = 0.1001513154359418
$ws.Range("J5").Value = 0.117124939437743
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.5776695000000001
$ws.Range("N5").Value = 1.155339
$ws.Range("O5").Value = 0.3527828075714595
$ws.Range("P5").Value = 0.2775942642546679
$ws.Range("Q5").Value = 0.1074829201785
$ws.Range("R5").Value = 0.6448975210710001
$ws.Range("S5").Value = 0.03533166224146639
$ws.Range("T5").Value = 0.03251321138909281
$ws.Range("I6").Value = 0.1001513154359418
$ws.Range("J6").Value = 0.117124939437743
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.04119366666666666
$ws.Range("N6").Value = 0.123581
$ws.Range("O6").Value = 0.02515697536537676
$ws.Range("P6").Value = 0.02969290984798064
$ws.Range("Q6").Value = 0.007664617201
$ws.Range("R6").Value = 0.06898155480900001
$ws.Range("S6").Value = 0.002519504175232065
$ws.Range("T6").Value = 0.003477780267675097
$ws.Range("I7").Value = 0.1001513154359418
$ws.Range("J7").Value = 0.117124939437743
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.1108493333333333
$ws.Range("N7").Value = 0.332548
$ws.Range("O7").Value = 0.06769569629478085
$ws.Range("P7").Value = 0.07990158506668717
$ws.Range("Q7").Value = 0.020624959508
$ws.Range("R7").Value = 0.185624635572
$ws.Range("S7").Value = 0.006779813033274312
$ws.Range("T7").Value = 0.009358468311915408
$ws.Range("G8").Value = 0.8076995
$ws.Range("H8").Value = 1.615399
$ws.Range("I8").Value = 0.4347568694579387
$ws.Range("J8").Value = 0.338959581867057
$ws.Range("M8").Value = 0.1727555
$ws.Range("N8").Value = 0.345511
$ws.Range("O8").Value = 0.1055017969849737
$ws.Range("P8").Value = 0.08301621587853829
$ws.Range("Q8").Value = 0.13953453097225
$ws.Range("R8").Value = 0.5581381238890001
$ws.Range("S8").Value = 0.04586763097937417
$ws.Range("T8").Value = 0.02813914182237467
$ws.Range("G9").Value = 0.8076995
$ws.Range("H9").Value = 1.615399
$ws.Range("I9").Value = 0.4347568694579387
$ws.Range("J9").Value = 0.338959581867057
$ws.Range("O9").Value = 0.2846509696390457
$ws.Range("P9").Value = 0.3359750310550052
$ws.Range("Q9").Value = 0.3764735831470001
$ws.Range("R9").Value = 2.258841498882
$ws.Range("S9").Value = 0.1237539644484383
$ws.Range("T9").Value = 0.113881956044176
$ws.Range("G10").Value = 0.8076995
$ws.Range("H10").Value = 1.615399
$ws.Range("I10").Value = 0.4347568694579387
$ws.Range("J10").Value = 0.338959581867057
$ws.Range("M10").Value = 0.268891
$ws.Range("N10").Value = 0.806673
$ws.Range("O10").Value = 0.1642117541443634
$ws.Range("P10").Value = 0.1938199938971208
$ws.Range("Q10").Value = 0.2171831262545
$ws.Range("R10").Value = 1.303098757527
$ws.Range("S10").Value = 0.07139218816000011
$ws.Range("T10").Value = 0.06569714408884361
$ws.Range("G11").Value = 0.8076995
$ws.Range("H11").Value = 1.615399
$ws.Range("I11").Value = 0.4347568694579387
$ws.Range("J11").Value = 0.338959581867057
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.5776695000000001
$ws.Range("N11").Value = 1.155339
$ws.Range("O11").Value = 0.3527828075714595
$ws.Range("P11").Value = 0.2775942642546679
$ws.Range("Q11").Value = 0.46658336631525
$ws.Range("R11").Value = 1.866333465261
$ws.Range("S11").Value = 0.1533747490183501
$ws.Range("T11").Value = 0.09409323574045556
$ws.Range("G12").Value = 0.8076995
$ws.Range("H12").Value = 1.615399
$ws.Range("I12").Value = 0.4347568694579387
$ws.Range("J12").Value = 0.338959581867057
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.04119366666666666
$ws.Range("N12").Value = 0.123581
$ws.Range("O12").Value = 0.02515697536537676
$ws.Range("P12").Value = 0.02969290984798064
$ws.Range("Q12").Value = 0.03327210396983333
$ws.Range("R12").Value = 0.199632623819
$ws.Range("S12").Value = 0.01093716785488168
$ws.Range("T12").Value = 0.01006469630648774
$ws.Range("G13").Value = 0.8076995
$ws.Range("H13").Value = 1.615399
$ws.Range("I13").Value = 0.4347568694579387
$ws.Range("J13").Value = 0.338959581867057
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.1108493333333333
$ws.Range("N13").Value = 0.332548
$ws.Range("O13").Value = 0.06769569629478085
$ws.Range("P13").Value = 0.07990158506668717
$ws.Range("Q13").Value = 0.08953295110866667
$ws.Range("R13").Value = 0.5371977066520001
$ws.Range("S13").Value = 0.0294311689968943
$ws.Range("T13").Value = 0.02708340786471937
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.8640563333333334
$ws.Range("H14").Value = 2.592169
$ws.Range("I14").Value = 0.4650918151061195
$ws.Range("J14").Value = 0.5439154786951998
$ws.Range("M14").Value = 0.1727555
$ws.Range("N14").Value = 0.345511
$ws.Range("O14").Value = 0.1055017969849737
$ws.Range("P14").Value = 0.08301621587853829
$ws.Range("Q14").Value = 0.1492704838931667
$ws.Range("R14").Value = 0.8956229033590001
$ws.Range("S14").Value = 0.04906802225669876
$ws.Range("T14").Value = 0.0451538047990392
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.8640563333333334
$ws.Range("H15").Value = 2.592169
$ws.Range("I15").Value = 0.4650918151061195
$ws.Range("J15").Value = 0.5439154786951998
$ws.Range("O15").Value = 0.2846509696390457
$ws.Range("P15").Value = 0.3359750310550052
$ws.Range("Q15").Value = 0.4027418413046668
$ws.Range("R15").Value = 3.624676571742001
$ws.Range("S15").Value = 0.1323888361411407
$ws.Range("T15").Value = 0.1827420198459178
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.8640563333333334
$ws.Range("H16").Value = 2.592169
$ws.Range("I16").Value = 0.4650918151061195
$ws.Range("J16").Value = 0.5439154786951998
$ws.Range("M16").Value = 0.268891
$ws.Range("N16").Value = 0.806673
$ws.Range("O16").Value = 0.1642117541443634
$ws.Range("P16").Value = 0.1938199938971208
$ws.Range("Q16").Value = 0.2323369715263333
$ws.Range("R16").Value = 2.091032743737
$ws.Range("S16").Value = 0.0763735427967618
$ws.Range("T16").Value = 0.1054216947612532
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.8640563333333334
$ws.Range("H17").Value = 2.592169
$ws.Range("I17").Value = 0.4650918151061195
$ws.Range("J17").Value = 0.5439154786951998
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.5776695000000001
$ws.Range("N17").Value = 1.155339
$ws.Range("O17").Value = 0.3527828075714595
$ws.Range("P17").Value = 0.2775942642546679
$ws.Range("Q17").Value = 0.4991389900485001
$ws.Range("R17").Value = 2.994833940291
$ws.Range("S17").Value = 0.164076396311643
$ws.Range("T17").Value = 0.1509878171251195
$ws.Range("E18").Value = 3
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 0.8640563333333334
$ws.Range("H18").Value = 2.592169
$ws.Range("I18").Value = 0.4650918151061195
$ws.Range("J18").Value = 0.5439154786951998
$ws.Range("K18").Value = 1
$ws.Range("L18").Value = 0.3333333333333333
$ws.Range("M18").Value = 0.04119366666666666
$ws.Range("N18").Value = 0.123581
$ws.Range("O18").Value = 0.02515697536537676
$ws.Range("P18").Value = 0.02969290984798064
$ws.Range("Q18").Value = 0.03559364857655556
$ws.Range("R18").Value = 0.320342837189
$ws.Range("S18").Value = 0.01170030333526301
$ws.Range("T18").Value = 0.0161504332738178
$ws.Range("E19").Value = 3
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 0.8640563333333334
$ws.Range("H19").Value = 2.592169
$ws.Range("I19").Value = 0.4650918151061195
$ws.Range("J19").Value = 0.5439154786951998
$ws.Range("K19").Value = 2
$ws.Range("L19").Value = 0.6666666666666666
$ws.Range("M19").Value = 0.1108493333333333
$ws.Range("N19").Value = 0.332548
$ws.Range("O19").Value = 0.06769569629478085
$ws.Range("P19").Value = 0.07990158506668717
$ws.Range("Q19").Value = 0.09578006851244446
$ws.Range("R19").Value = 0.8620206166120001
$ws.Range("S19").Value = 0.03148471426461224
$ws.Range("T19").Value = 0.04345970889005239

Write-Host "Updated 226 cells with new TPM values"
